$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the client's occupation (column F, "ocupacion") for the first
# data row from "estudiante" to "maestro".
$ws.Range("F2").Value = "maestro"
